$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 24.10000000000033
$ws.Range("G2").Value = [double]"1.110223024625157e-16"
$ws.Range("H2").Value = [double]"4.662353909187009e-16"
$ws.Range("K2").Value = 452.9059165541908
$ws.Range("L2").Value = "[361.6305809805874, 544.1812521277941]"
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 1.566079220708425
$ws.Range("P2").Value = "[1.3522370781217319, 1.7799213632951183]"
$ws.Range("S2").Value = 823.1329243291839
$ws.Range("T2").Value = "[762.0787323913034, 884.1871162670643]"
$ws.Range("W2").Value = 18.09309309309334
$ws.Range("X2").Value = 17.27287287287311
$ws.Range("Y2").Value = 18.91331331331357

# Row 3
$ws.Range("C3").Value = "2_induction_dd"
$ws.Range("E3").Value = 23.85000000000029
$ws.Range("H3").Value = [double]"4.662353909187009e-16"
$ws.Range("I3").Value = 0.001266727310926918
$ws.Range("K3").Value = 253.8204589722343
$ws.Range("L3").Value = "[213.64159811621641, 293.99931982825217]"
$ws.Range("O3").Value = 2.270500396288118
$ws.Range("P3").Value = "[2.11955300152104, 2.421447791055196]"
$ws.Range("S3").Value = 707.0312904389067
$ws.Range("T3").Value = "[686.6019064862714, 727.4606743915419]"
$ws.Range("W3").Value = 15.23153153153172
$ws.Range("X3").Value = 14.65855855855873
$ws.Range("Y3").Value = 15.8045045045047

# Row 4
$ws.Range("C4").Value = "3_hypo_dd"
$ws.Range("E4").Value = 23.70000000000027
$ws.Range("G4").Value = [double]"1.159506934911292e-10"
$ws.Range("H4").Value = [double]"2.806673367410086e-10"
$ws.Range("I4").Value = [double]"1.110223024625157e-16"
$ws.Range("K4").Value = 129.6573846840802
$ws.Range("L4").Value = "[85.94007736177969, 173.37469200638066]"
$ws.Range("M4").Value = [double]"1.035688845796301e-08"
$ws.Range("N4").Value = [double]"1.035688845796301e-08"
$ws.Range("O4").Value = 1.314500229429964
$ws.Range("P4").Value = "[0.9119738433844251, 1.7170266154755023]"
$ws.Range("Q4").Value = [double]"3.360474121194557e-10"
$ws.Range("R4").Value = [double]"3.360474121194557e-10"
$ws.Range("S4").Value = 659.9581717835814
$ws.Range("T4").Value = "[634.0343368885558, 685.882006678607]"
$ws.Range("W4").Value = 18.74174174174195
$ws.Range("X4").Value = 17.22342342342362
$ws.Range("Y4").Value = 20.26006006006028

# Row 5
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "4_hypo_dl"
$ws.Range("E5").Value = 25.52000000000055
$ws.Range("G5").Value = [double]"1.110223024625157e-16"
$ws.Range("H5").Value = [double]"4.662353909187009e-16"
$ws.Range("K5").Value = 308.7776886554807
$ws.Range("L5").Value = "[256.9849643659252, 360.5704129450362]"
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = -2.138421425866927
$ws.Range("P5").Value = "[-2.3145267197618495, -1.962316131972004]"
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 652.5730014004326
$ws.Range("T5").Value = "[625.6545208149248, 679.4914819859405]"
$ws.Range("W5").Value = 8.685485485485675
$ws.Range("X5").Value = 7.970210210210386
$ws.Range("Y5").Value = 9.400760760760964

# Remove row 6 (merged into the dataset change / dimension shrink)
$ws.Rows.Item(6).Delete()

Write-Output "done"